$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("reviews_count") entirely; subsequent columns (F..K) shift left to (E..J)
$ws.Range("E1").EntireColumn.Delete()
